$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the now-obsolete rows (old rows 8, 9, 10 drop out of the feed).
#    Do this first so the remaining rows keep their natural row numbers.
# ---------------------------------------------------------------------------
$ws.Rows("8:10").Delete()

# ---------------------------------------------------------------------------
# 2. Refresh the data for rows 2-7 with the newly scraped opportunities.
# ---------------------------------------------------------------------------
$ids         = @("1331894","1331888","1331887","1331869","1331584","1301473")
$links       = @(
  "https://aiesec.org/opportunity/global-talent/1331894",
  "https://aiesec.org/opportunity/global-talent/1331888",
  "https://aiesec.org/opportunity/global-talent/1331887",
  "https://aiesec.org/opportunity/global-talent/1331869",
  "https://aiesec.org/opportunity/global-talent/1331584",
  "https://aiesec.org/opportunity/global-talent/1301473"
)
$titles      = @(
  "Marketing",
  "Research Assistant",
  "Sales Specialist",
  "Language Specialist - Spanish",
  "Web Developer Intern",
  "Copy Writer in Istanbul"
)
$countries   = @(
  "Kartepe, Kocaeli, Türkiye",
  "Hong Kong",
  "Başakşehir, Başak, 34490 Başakşehir/İstanbul, Türkiye",
  "Colombo, Sri Lanka",
  "Phagwara, Punjab, India",
  "İstanbul, Türkiye"
)
$premiums    = @("No","No","No","No","No","No")
$applicants  = @("2 applicants","2 applicants","5 applicants","2 applicants","0 applicants","39 applicants")
$durations   = @("3 - 6 Months","6 - 18 Months","6 - 18 Months","6 - 18 Months","6 - 18 Months","6 - 18 Months")
$orgs        = @(
  "Atom Mek",
  "WSC Holding Limited",
  "HASEL AMB.MAK.VE KONVEYÖR SİS.SAN.VE TİCLTD.ŞTİ.",
  "Aitken Spence Travels (Pvt) Ltd",
  "GNA University",
  "Raff Textile"
)

for ($i = 0; $i -lt 6; $i++) {
  $r = $i + 2

  # Column A holds an opportunity id that looks numeric ("1331894"); force it
  # to stay text (matching the rest of the id column) instead of silently
  # becoming a number, then drop the temporary number format again.
  $idCell = $ws.Cells.Item($r, 1)
  $idCell.NumberFormat = "@"
  $idCell.Value = $ids[$i]
  $idCell.ClearFormats()

  $ws.Cells.Item($r, 2).Value = $links[$i]
  $ws.Cells.Item($r, 3).Value = $titles[$i]
  $ws.Cells.Item($r, 4).Value = $countries[$i]
  $ws.Cells.Item($r, 5).Value = $premiums[$i]
  $ws.Cells.Item($r, 6).Value = $applicants[$i]
  $ws.Cells.Item($r, 7).Value = $durations[$i]
  $ws.Cells.Item($r, 8).Value = $orgs[$i]
}

# ---------------------------------------------------------------------------
# 3. The "Yes" premium flag (and its yellow highlight) on the old row 2 is
#    gone now that the row holds a "No" opportunity - strip the highlight.
# ---------------------------------------------------------------------------
$ws.Range("E2").ClearFormats()

# ---------------------------------------------------------------------------
# 4. A couple of columns were resized to fit the new copy.
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 31.1666666667
$ws.Columns("D").ColumnWidth = 55.1666666667
$ws.Columns("H").ColumnWidth = 50.1666666667
